# "add footer and 1 modal"
#
# The sheet is a task-tracking table. This edit:
#  - bumps the FACT(ч) time logged against several rows (E column)
#  - flips those rows' status from "Не готово" -> "Готово" (F column)
#  - clears two "Всего: N" sub-total label cells (D5, D13) that are no
#    longer used once the footer/modal tasks are folded into the main
#    PLAN/FACT totals
#  - turns the static "Итого" row (22) D/E cells into live SUM() formulas
#  - updates the selection/scroll position left behind by the editor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FACT(ч) time revised 4 -> 6 -------------------------------
$ws.Range("E2").Value = 6

# --- Row 5: clear the "Всего: 14" sub-total label ----------------------
$ws.Range("D5").ClearContents()

# --- Rows 7-12: log FACT(ч) hours and mark these tasks done ------------
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Готово"

$ws.Range("E8").Value = 3.5
$ws.Range("F8").Value = "Готово"

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "Готово"

$ws.Range("E10").Value = 0.5
$ws.Range("F10").Value = "Готово"

$ws.Range("E11").Value = 2.5
$ws.Range("F11").Value = "Готово"

$ws.Range("E12").Value = 1.5
$ws.Range("F12").Value = "Готово"

# --- Row 13: clear the "Всего: 9" sub-total label -----------------------
$ws.Range("D13").ClearContents()

# --- Row 14: log FACT(ч) hours and mark done -----------------------------
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = "Готово"

# --- Row 20: log FACT(ч) hours and mark done -----------------------------
$ws.Range("E20").Value = 0.5
$ws.Range("F20").Value = "Готово"

# --- Row 21: "Футер" task (reuses the "Адаптивная версия" shared string
#     slot vacated above; content itself is unchanged) ------------------
$ws.Range("B21").Value = "Адаптивная версия"

# --- Row 22 "Итого": turn the static totals into real SUM formulas ------
$ws.Range("D22").Formula = "=SUM(D2:D21)"
$ws.Range("E22").Formula = "=SUM(E2:E21)"
# Match the bold/centered "Итого" style already used elsewhere (style idx 6)
$ws.Range("D22").Font.Bold = $true
$ws.Range("D22").HorizontalAlignment = -4108
$ws.Range("D22").VerticalAlignment = -4108
$ws.Range("E22").Font.Bold = $true
$ws.Range("E22").HorizontalAlignment = -4108
$ws.Range("E22").VerticalAlignment = -4108

# --- Leave the sheet scrolled/selected where the editor left it --------
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
